$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point: insert two rows (Primera + Segunda) right after the
# existing row 50, pushing all subsequent rows down by two (old row 51 ->
# new row 53, ..., old row 75 -> new row 77).
$ws.Rows("51:52").Insert()

# Row 51 - "Primera" quality entry for 2022-01-11
$ws.Cells.Item(51, 1).Value = 11
$ws.Cells.Item(51, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(51, 3).Value = "Bíobío"
$ws.Cells.Item(51, 4).Value = [DateTime]"2022-01-11"
$ws.Cells.Item(51, 5).Value = 8
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100101
$ws.Cells.Item(51, 8).Value = "Berries"
$ws.Cells.Item(51, 9).Value = 100101001
$ws.Cells.Item(51, 10).Value = "Arándano (blue)"
$ws.Cells.Item(51, 11).Value = "Sin especificar"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 200
$ws.Cells.Item(51, 14).Value = 3500
$ws.Cells.Item(51, 15).Value = 4000
$ws.Cells.Item(51, 16).Value = 3750
$ws.Cells.Item(51, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(51, 18).Value = "Región de Ñuble"
$ws.Cells.Item(51, 19).Value = 1875
$ws.Cells.Item(51, 20).Value = 2

# Row 52 - "Segunda" quality entry for 2022-01-11
$ws.Cells.Item(52, 1).Value = 11
$ws.Cells.Item(52, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(52, 3).Value = "Bíobío"
$ws.Cells.Item(52, 4).Value = [DateTime]"2022-01-11"
$ws.Cells.Item(52, 5).Value = 8
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100101
$ws.Cells.Item(52, 8).Value = "Berries"
$ws.Cells.Item(52, 9).Value = 100101001
$ws.Cells.Item(52, 10).Value = "Arándano (blue)"
$ws.Cells.Item(52, 11).Value = "Sin especificar"
$ws.Cells.Item(52, 12).Value = "Segunda"
$ws.Cells.Item(52, 13).Value = 100
$ws.Cells.Item(52, 14).Value = 3000
$ws.Cells.Item(52, 15).Value = 3000
$ws.Cells.Item(52, 16).Value = 3000
$ws.Cells.Item(52, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(52, 18).Value = "Región de Ñuble"
$ws.Cells.Item(52, 19).Value = 1500
$ws.Cells.Item(52, 20).Value = 2
